# Add data for 2022-04-30:
# - Advance the "through" date referenced in the sheet name and the
#   April row label from 04-21 to 04-22.
# - Update April's monthly counts (row 5) and the Total row (row 6)
#   to reflect the newly added day's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-04-22"

# Update the April row label.
$ws.Range("A5").Value = "April (through 04-22)"

# Update April row (row 5) year-by-year counts.
$ws.Range("B5").Value = 13
$ws.Range("C5").Value = 23
$ws.Range("D5").Value = 43
$ws.Range("E5").Value = 40
$ws.Range("F5").Value = 33
$ws.Range("G5").Value = 44
$ws.Range("H5").Value = 82
$ws.Range("I5").Value = 92

# Update Total row (row 6) year-by-year counts.
$ws.Range("B6").Value = 79
$ws.Range("C6").Value = 151
$ws.Range("D6").Value = 232
$ws.Range("E6").Value = 237
$ws.Range("F6").Value = 143
$ws.Range("G6").Value = 242
$ws.Range("H6").Value = 505
$ws.Range("I6").Value = 528
